# Iteration plan workbook update:
#  - Use Cases sheet: add UC_2 (info collection/management) and UC_3 (data
#    analysis/visualization) use-case rows.
#  - Sprint 1 sheet: fill in Estimated/Real Effort numbers for the
#    single-question operations (moveUp / moveDown / remove) tasks that were
#    previously left blank, now that that work (and the "complex modal" test)
#    has effort estimates.
#  - Refresh a few view-state bits (zoom, selection) left behind by the
#    authoring session.

$wb = $excel.ActiveWorkbook

# --- "Use Cases" sheet: append UC_2 and UC_3 --------------------------------
$wsUseCases = $wb.Worksheets.Item("Use Cases")

$wsUseCases.Cells.Item(3, 1).Value = "UC_2"
$wsUseCases.Cells.Item(3, 2).Value = "信息收集和管理"
$wsUseCases.Cells.Item(3, 3).Value = "活动发布者可以通过链接发布问卷，之后可以查看、查询、筛选，并导出数据。活动参与者通过链接填写问卷，后台数据库实时更新"

$wsUseCases.Cells.Item(4, 1).Value = "UC_3"
$wsUseCases.Cells.Item(4, 2).Value = "数据整理、分析、可视化"
$wsUseCases.Cells.Item(4, 3).Value = "活动发布者可以查看不同ip的输入信息，并对数据进行简单的分析，可以导出原始数据和分析结果"

# Description column needed more room for the new, longer use-case text.
$wsUseCases.Columns.Item(3).ColumnWidth = 110.85714285714286

# --- "Sprint 1" sheet: fill in effort estimates for single-question ops ----
$wsSprint1 = $wb.Worksheets.Item("Sprint 1")

$wsSprint1.Cells.Item(7, 4).Value = 2        # moveUp/moveDown/remove modal test
$wsSprint1.Cells.Item(8, 4).Value = 0.5
$wsSprint1.Cells.Item(11, 4).Value = 1
$wsSprint1.Cells.Item(14, 4).Value = 1
$wsSprint1.Cells.Item(15, 5).Value = 1
$wsSprint1.Cells.Item(20, 4).Value = 1.5

# --- View-state refresh ------------------------------------------------------
$wsUseCases.Activate()
$excel.ActiveWindow.Zoom = 125
$wsUseCases.Range("C4").Select() | Out-Null

$wsProductBacklog = $wb.Worksheets.Item("Product Backlog")
$wsProductBacklog.Activate()
$excel.ActiveWindow.Zoom = 125

$wsSprint1.Activate()
$excel.ActiveWindow.Zoom = 125
$wsSprint1.Range("E31").Select() | Out-Null
